# Update the JULY_2024 attendance sheet: a second class date was added,
# doubling "Total Classes" (G12) from 3 to 6, and recording attendance
# for that second date in column E for each student row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JULY_2024")

# Header block: number of classes held on the newly recorded date (row 10)
$ws.Range("E10").Value = 14

# Summary "Total Classes" row: second-date class count (E12) and the new
# grand total of classes held this month (G12: 3 -> 6)
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 6

# Per-student attendance for the second class date (column E), rows 14-24.
# Columns G (Total classes attended) and H (%age) are formulas that will
# recalculate automatically once the underlying divisor (row 12 total,
# referenced literally inside each H-column formula) is updated below.
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3
$ws.Range("E22").Value = 3
$ws.Range("E23").Value = 0
$ws.Range("E24").Value = 3

# The %age formulas divide by the old fixed total of 3 classes; repoint
# them at the new total of 6 classes now that a second date was added.
$ws.Range("H14").Formula = "=(G14/6)*100"
$ws.Range("H15").Formula = "=(G15/6)*100"
$ws.Range("H16").Formula = "=(G16/6)*100"
$ws.Range("H17").Formula = "=(G17/6)*100"
$ws.Range("H18").Formula = "=(G18/6)*100"
$ws.Range("H19").Formula = "=(G19/6)*100"
$ws.Range("H20").Formula = "=(G20/6)*100"
$ws.Range("H21").Formula = "=(G21/6)*100"
$ws.Range("H22").Formula = "=(G22/6)*100"
$ws.Range("H23").Formula = "=(G23/6)*100"
$ws.Range("H24").Formula = "=(G24/6)*100"

# Selection/view tweaks to match the author's final on-screen state.
$ws.Range("H20").Select()

$wb.Save()
